# modified modules watershed 4
#
# The "wtrsmple_de"/sampling-period column (A) used to hold an actual date
# serial (formatted as a date). It is now recorded simply as the text year
# "2017" for every data row, so replace the date values in A2:A3 with the
# literal text "2017" (not a number) and drop the now-unused date number
# format from those two cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as text instead of
# re-parsing "2017" back into a number.
$ws.Range("A2").Value = "'2017"
$ws.Range("A3").Value = "'2017"

# Drop the old date-format styling on just these two cells (leaves the
# rest of the rows' number formats, e.g. lat/long in E:G, untouched).
$ws.Range("A2").ClearFormats()
$ws.Range("A3").ClearFormats()
